$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value into a cell while keeping it stored
# as text (matching the workbook's original inlineStr cell type for column D),
# and then reset the cell style so no stray number-format/quote-prefix style
# is left behind on the cell.
function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Column D (Price) updates
Set-TextValue $ws.Range("D2")  "228.84"
Set-TextValue $ws.Range("D3")  "22.42"
Set-TextValue $ws.Range("D4")  "5.292"
Set-TextValue $ws.Range("D5")  "0.05538"
Set-TextValue $ws.Range("D6")  "3.385"
Set-TextValue $ws.Range("D7")  "6.473"
Set-TextValue $ws.Range("D8")  "1.068"
Set-TextValue $ws.Range("D9")  "0.7706"
Set-TextValue $ws.Range("D10") "0.1378"
Set-TextValue $ws.Range("D11") "0.07425"
Set-TextValue $ws.Range("D13") "0.02949"
Set-TextValue $ws.Range("D14") "0.09263"
Set-TextValue $ws.Range("D15") "0.001648"
Set-TextValue $ws.Range("D16") "3.257"
Set-TextValue $ws.Range("D17") "0.04774"
Set-TextValue $ws.Range("D18") "0.0005895"
Set-TextValue $ws.Range("D19") "0.006225"
Set-TextValue $ws.Range("D20") "0.005227"
Set-TextValue $ws.Range("D23") "3.885"
Set-TextValue $ws.Range("D26") "0.1286"
Set-TextValue $ws.Range("D27") "0.0005004"

$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"

Set-TextValue $ws.Range("D40") "0.03949"
Set-TextValue $ws.Range("D41") "0.007119"

$ws.Range("E41").Value = "40KickTokenKICK"

# Rows 42/43: CEJI and BKEXToken swap places
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1035"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.002770"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"

Set-TextValue $ws.Range("D44") "0.008744"
Set-TextValue $ws.Range("D45") "0.00005446"
Set-TextValue $ws.Range("D46") "0.00000000751"
Set-TextValue $ws.Range("D47") "0.7859"
Set-TextValue $ws.Range("D48") "0.08787"

$ws.Range("E48").Value = "47BOLOBOLO"

Set-TextValue $ws.Range("D49") "0.00002102"
Set-TextValue $ws.Range("D50") "0.01011"
